$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells that receive numeric-looking text to remain text,
# matching the source workbook where these are stored as inline strings.
# (Applied as separate statements per contiguous block; union ranges with commas
# are not reliably supported for NumberFormat assignment.)
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D8:D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31:D34").NumberFormat = "@"
$ws.Range("D36:D42").NumberFormat = "@"
$ws.Range("D44:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.397.41"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "1.872.23"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "244.57"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("D8").Value = "0.2875"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "0.06490"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").Value = "21.87"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "99.73"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "0.07788"
$ws.Range("D13").Value = "1.872.26"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "0.7308"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "5.171"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "286.23"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "30.383.91"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "13.12"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "0.000007492"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "2.114.71"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "5.313"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "6.322"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "163.08"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "9.050"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").Value = "18.92"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "0.09666"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").Value = "1.486"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "4.231"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").Value = "4.147"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").Value = "0.04805"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "0.6887"
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("D37").Value = "2.726"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "0.01900"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").Value = "2.843"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").Value = "75.90"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").Value = "6.278"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "1.964"
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "0.9988"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "0.8238"
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("D46").Value = "101.16"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "9.813"
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("D48").Value = "7.014"
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D49").Value = "35.01"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").Value = "0.05765"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "886.09"
$ws.Range("E51").Value = "  -4.13%  "
